$wb = $excel.ActiveWorkbook

# Sheet1 = 展览 (Exhibitions), Sheet2 = 演出 (Performances),
# Sheet3 = 本地生活 (Local life), Sheet4 = 全部类型 (All types - combines the above)
$sheets = @{}
$sheets[1] = $wb.Worksheets.Item(1)
$sheets[2] = $wb.Worksheets.Item(2)
$sheets[3] = $wb.Worksheets.Item(3)
$sheets[4] = $wb.Worksheets.Item(4)

# Update column F (想去人数 / "want to go" count) values per the refreshed data snapshot.
$sheets[1].Cells.Item(2, 6).Value = 302
$sheets[1].Cells.Item(3, 6).Value = 484
$sheets[1].Cells.Item(4, 6).Value = 219
$sheets[1].Cells.Item(5, 6).Value = 244
$sheets[1].Cells.Item(6, 6).Value = 294
$sheets[1].Cells.Item(7, 6).Value = 7431
$sheets[1].Cells.Item(8, 6).Value = 86
$sheets[1].Cells.Item(10, 6).Value = 3584
$sheets[1].Cells.Item(12, 6).Value = 580
$sheets[1].Cells.Item(13, 6).Value = 594
$sheets[1].Cells.Item(14, 6).Value = 434
$sheets[1].Cells.Item(15, 6).Value = 131
$sheets[1].Cells.Item(16, 6).Value = 89
$sheets[1].Cells.Item(17, 6).Value = 748
$sheets[1].Cells.Item(18, 6).Value = 19
$sheets[1].Cells.Item(19, 6).Value = 55
$sheets[1].Cells.Item(20, 6).Value = 213
$sheets[1].Cells.Item(22, 6).Value = 228
$sheets[1].Cells.Item(23, 6).Value = 125
$sheets[1].Cells.Item(24, 6).Value = 370
$sheets[1].Cells.Item(25, 6).Value = 130
$sheets[1].Cells.Item(26, 6).Value = 1072
$sheets[1].Cells.Item(27, 6).Value = 74
$sheets[1].Cells.Item(28, 6).Value = 191
$sheets[1].Cells.Item(29, 6).Value = 2117
$sheets[1].Cells.Item(30, 6).Value = 639
$sheets[1].Cells.Item(31, 6).Value = 26
$sheets[1].Cells.Item(32, 6).Value = 27
$sheets[1].Cells.Item(34, 6).Value = 580
$sheets[1].Cells.Item(35, 6).Value = 31

$sheets[2].Cells.Item(3, 6).Value = 62
$sheets[2].Cells.Item(6, 6).Value = 74

$sheets[3].Cells.Item(2, 6).Value = 410

$sheets[4].Cells.Item(2, 6).Value = 410
$sheets[4].Cells.Item(3, 6).Value = 302
$sheets[4].Cells.Item(4, 6).Value = 484
$sheets[4].Cells.Item(5, 6).Value = 219
$sheets[4].Cells.Item(6, 6).Value = 244
$sheets[4].Cells.Item(7, 6).Value = 294
$sheets[4].Cells.Item(8, 6).Value = 7431
$sheets[4].Cells.Item(9, 6).Value = 86
$sheets[4].Cells.Item(12, 6).Value = 3585
$sheets[4].Cells.Item(14, 6).Value = 580
$sheets[4].Cells.Item(15, 6).Value = 594
$sheets[4].Cells.Item(16, 6).Value = 434
$sheets[4].Cells.Item(17, 6).Value = 62
$sheets[4].Cells.Item(18, 6).Value = 131
$sheets[4].Cells.Item(19, 6).Value = 90
$sheets[4].Cells.Item(22, 6).Value = 74
$sheets[4].Cells.Item(23, 6).Value = 748
$sheets[4].Cells.Item(24, 6).Value = 19
$sheets[4].Cells.Item(25, 6).Value = 55
$sheets[4].Cells.Item(26, 6).Value = 213
$sheets[4].Cells.Item(31, 6).Value = 228
$sheets[4].Cells.Item(32, 6).Value = 125
$sheets[4].Cells.Item(33, 6).Value = 370
$sheets[4].Cells.Item(34, 6).Value = 130
$sheets[4].Cells.Item(35, 6).Value = 1072
$sheets[4].Cells.Item(36, 6).Value = 74
$sheets[4].Cells.Item(37, 6).Value = 191
$sheets[4].Cells.Item(38, 6).Value = 2117
$sheets[4].Cells.Item(39, 6).Value = 639
$sheets[4].Cells.Item(40, 6).Value = 26
$sheets[4].Cells.Item(41, 6).Value = 27
$sheets[4].Cells.Item(43, 6).Value = 580
$sheets[4].Cells.Item(44, 6).Value = 31
